$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1863905325443787
$ws.Range("C2").Value = 0.5710059171597633
$ws.Range("J2").Value = 0.01183431952662722
$ws.Range("P2").Value = 0.1420118343195266
$ws.Range("S2").Value = 0.08875739644970414
$ws.Range("C3").Value = 0.03980099502487562
$ws.Range("J3").Value = 0.03980099502487562
$ws.Range("P3").Value = 0.746268656716418
$ws.Range("S3").Value = 0.1741293532338309
$ws.Range("B6").Value = 0.05633802816901409
$ws.Range("D6").Value = 0.02816901408450704
$ws.Range("F6").Value = 0.06103286384976526
$ws.Range("J6").Value = 0.3051643192488263
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.1408450704225352
$ws.Range("R6").Value = 0.06572769953051644
$ws.Range("S6").Value = 0.3286384976525822
$ws.Range("B7").Value = 0.1276595744680851
$ws.Range("D7").Value = 0.02553191489361702
$ws.Range("F7").Value = 0.05531914893617021
$ws.Range("J7").Value = 0.1276595744680851
$ws.Range("O7").Value = 0.03404255319148936
$ws.Range("Q7").Value = 0.1872340425531915
$ws.Range("R7").Value = 0.1319148936170213
$ws.Range("S7").Value = 0.3106382978723404
$ws.Range("B8").Value = 0.1002277904328018
$ws.Range("D8").Value = 0.02733485193621868
$ws.Range("E8").Value = 0.002277904328018223
$ws.Range("F8").Value = 0.05922551252847381
$ws.Range("J8").Value = 0.1093394077448747
$ws.Range("O8").Value = 0.01366742596810934
$ws.Range("Q8").Value = 0.1981776765375854
$ws.Range("R8").Value = 0.10250569476082
$ws.Range("S8").Value = 0.387243735763098
$ws.Range("B9").Value = 0.08450704225352113
$ws.Range("D9").Value = 0.009389671361502348
$ws.Range("E9").Value = 0.009389671361502348
$ws.Range("F9").Value = 0.08450704225352113
$ws.Range("J9").Value = 0.1220657276995305
$ws.Range("O9").Value = 0.009389671361502348
$ws.Range("Q9").Value = 0.1643192488262911
$ws.Range("R9").Value = 0.1032863849765258
$ws.Range("S9").Value = 0.4131455399061033
$ws.Range("B10").Value = 0.1286821705426357
$ws.Range("D10").Value = 0.01395348837209302
$ws.Range("E10").Value = 0.001550387596899225
$ws.Range("F10").Value = 0.06124031007751938
$ws.Range("J10").Value = 0.1131782945736434
$ws.Range("O10").Value = 0.01162790697674419
$ws.Range("Q10").Value = 0.2015503875968992
$ws.Range("R10").Value = 0.110077519379845
$ws.Range("S10").Value = 0.3581395348837209
$ws.Range("G11").Value = 0.1454545454545454
$ws.Range("J11").Value = 0.06363636363636363
$ws.Range("K11").Value = 0.1757575757575758
$ws.Range("L11").Value = 0.603030303030303
$ws.Range("S11").Value = 0.01212121212121212
$ws.Range("G12").Value = 0.7688442211055276
$ws.Range("J12").Value = 0.1909547738693467
$ws.Range("K12").Value = 0.01005025125628141
$ws.Range("L12").Value = 0.01005025125628141
$ws.Range("S12").Value = 0.02010050251256281
$ws.Range("G13").Value = 0.7678571428571429
$ws.Range("J13").Value = 0.2321428571428572
$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1422018348623853
$ws.Range("I15").Value = 0.0963302752293578
$ws.Range("J15").Value = 0.3807339449541284
$ws.Range("K15").Value = 0.06422018348623854
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("O15").Value = 0.03669724770642202
$ws.Range("S15").Value = 0.2431192660550459
$ws.Range("F16").Value = 0.02678571428571428
$ws.Range("H16").Value = 0.15625
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4419642857142857
$ws.Range("K16").Value = 0.1339285714285714
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("O16").Value = 0.04464285714285714
$ws.Range("S16").Value = 0.1071428571428571
$ws.Range("F17").Value = 0.01348314606741573
$ws.Range("H17").Value = 0.150561797752809
$ws.Range("I17").Value = 0.09662921348314607
$ws.Range("J17").Value = 0.4044943820224719
$ws.Range("K17").Value = 0.1393258426966292
$ws.Range("M17").Value = 0.01573033707865169
$ws.Range("O17").Value = 0.06292134831460675
$ws.Range("S17").Value = 0.1168539325842697
$ws.Range("F18").Value = 0.0158102766798419
$ws.Range("H18").Value = 0.1936758893280632
$ws.Range("I18").Value = 0.08300395256916997
$ws.Range("J18").Value = 0.4150197628458498
$ws.Range("K18").Value = 0.1383399209486166
$ws.Range("M18").Value = 0.01976284584980237
$ws.Range("O18").Value = 0.03557312252964427
$ws.Range("S18").Value = 0.09881422924901186
$ws.Range("F19").Value = 0.0199501246882793
$ws.Range("H19").Value = 0.2152950955943475
$ws.Range("I19").Value = 0.09393183707398171
$ws.Range("J19").Value = 0.3599334995843724
$ws.Range("K19").Value = 0.1039068994181214
$ws.Range("M19").Value = 0.03158769742310889
$ws.Range("O19").Value = 0.0798004987531172
$ws.Range("S19").Value = 0.09559434746467166
